$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column E (shifts E.. onward to F.. etc.)
$ws.Columns("E").Insert()

# Populate the new column E header rows (2-5) with the new field info
$ws.Range("E2").Value = "Resource Identifying Information"
$ws.Range("E3").Value = "Resource URI REQUIRED IF NO EAD ID"
$ws.Range("E4").Value = "res_uri"
$ws.Range("E5").Value = "Resource URI"

# Match styles used elsewhere for this header block
$ws.Range("F2").Copy()
$ws.Range("E2").PasteSpecial(-4122)
$ws.Range("F3").Copy()
$ws.Range("E3").PasteSpecial(-4122)
$ws.Range("N4").Copy()
$ws.Range("E4").PasteSpecial(-4122)
$ws.Range("F5").Copy()
$ws.Range("E5").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("E1:E1048576").Select()
